$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.853.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7666"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3126"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07122"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7625"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.925.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.361"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.128"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.958.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007817"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.016"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1623"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.389"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.033"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.510"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.540"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.504"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.118"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05438"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7447"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.701"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01946"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.781"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4463"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.100.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.079"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8554"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.868"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.657"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.046"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.040.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06084"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "
